# feat: add 2022-Q1 data
#
# - insert a new "2022-Q1" worksheet right after "2021-Q4" (same column
#   layout/formatting as "2021-Q4") and populate it with the quarter's
#   fund-holding rows
# - update the "总计" (totals) worksheet: push the existing "2021-Q4"
#   summary row down one row and add a new summary row for "2022-Q1"

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "2021-Q4"

# xlPasteFormats - used to copy over header / index-column styling only
$xlPasteFormats = -4122

function Set-TextCell($cell, [string]$text) {
    # Force text storage for numeric-looking strings (fund codes, ratios,
    # …) so leading zeros / exact decimal text survive, then drop the
    # scratch "@" number format again so the cell ends up unstyled, same
    # as the source data.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q1" sheet right after "2021-Q4"
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "2022-Q1"

# Copy header-row (B1:H1) and index-column (A) styling from "2021-Q4"
$ws1.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial($xlPasteFormats)
$ws1.Range("A2").Copy()
$newSheet.Range("A2:A6").PasteSpecial($xlPasteFormats)

$newSheet.Cells.Item(1, 2).Value = "基金代码"
$newSheet.Cells.Item(1, 3).Value = "基金名称"
$newSheet.Cells.Item(1, 4).Value = "基金规模"
$newSheet.Cells.Item(1, 5).Value = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value = "仓位占比"
$newSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value = "仓位排名"

# columns: 基金代码, 基金名称, 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名
$fundRows = @(
    @("519029", "华夏稳增混合",      "10.92", "92.99", "6.07", "0.6628", 2),
    @("009874", "九泰久睿量化股票",  "3.97",  "94.17", "3.17", "0.1258", 7),
    @("010120", "九泰久福量化股票A", "0.68",  "94.26", "3.18", "0.0216", 8),
    @("009043", "九泰久信量化股票",  "0.52",  "94.26", "3.17", "0.0165", 9),
    @("010121", "九泰久福量化股票C", "0.04",  "94.26", "3.18", "0.0013", 8)
)

for ($i = 0; $i -lt $fundRows.Count; $i++) {
    $row = $i + 2
    $r = $fundRows[$i]
    $newSheet.Cells.Item($row, 1).Value = $i
    Set-TextCell $newSheet.Cells.Item($row, 2) $r[0]
    Set-TextCell $newSheet.Cells.Item($row, 3) $r[1]
    Set-TextCell $newSheet.Cells.Item($row, 4) $r[2]
    Set-TextCell $newSheet.Cells.Item($row, 5) $r[3]
    Set-TextCell $newSheet.Cells.Item($row, 6) $r[4]
    Set-TextCell $newSheet.Cells.Item($row, 7) $r[5]
    $newSheet.Cells.Item($row, 8).Value = $r[6]
}

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: shift the existing "2021-Q4" row down one
#    row, then write the new "2022-Q1" summary row in its old spot
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

# move the old row 2 ("2021-Q4": 1, 0.48) down to row 3, preserving style
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial($xlPasteFormats)
$wsTotal.Cells.Item(3, 1).Value = 1
Set-TextCell $wsTotal.Cells.Item(3, 2) "2021-Q4"
$wsTotal.Cells.Item(3, 3).Value = 1
$wsTotal.Cells.Item(3, 4).Value = 0.48

# write the new "2022-Q1" summary row in row 2
Set-TextCell $wsTotal.Cells.Item(2, 2) "2022-Q1"
$wsTotal.Cells.Item(2, 3).Value = 5
$wsTotal.Cells.Item(2, 4).Value = 0.83

# ---------------------------------------------------------------------
# 3. Restore the original active sheet/selection ("2021-Q4") since the
#    source workbook's view state was untouched by this edit
# ---------------------------------------------------------------------
$ws1.Select()
